$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Price (D) and Volume(1h) (E) columns are stored as literal text in this
# sheet (e.g. "298.53", "-6.21%"), not as numbers/percentages. Force the
# target cells to Text format before assigning so Excel does not silently
# reinterpret the numeric- or percent-looking strings as real numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "299.32"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-6.23%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.13"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.33%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.970"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.98%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07908"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.87%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.906"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-11.91%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.020"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-2.59%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "7.724"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-4.16%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "3.85%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9227"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.67%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1095"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "8.60%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1818"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.46%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09180"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.18%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03533"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.80%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09869"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.62%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001396"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.53%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005685"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.44%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.489"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.79%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3441"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.99%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1308"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.70%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.069"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.20%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2399"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "8.85%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04496"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-2.54%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001212"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.59%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004589"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-3.27%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001250"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-4.00%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-6.90%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01882"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-3.93%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04677"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-6.06%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007568"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.12%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009552"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "25.48%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1321"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-5.64%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002119"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2.44%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01120"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.39%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006005"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-4.86%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.12%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-31.41%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.12%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001999"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.12%"
